$wb = $excel.ActiveWorkbook

# --- Workbook window view ---
$wb.Windows.Item(1).Width = 13920

# --- Sheet1 (BP4D): fill rows 20 and 21, rows 22-28 formatting, sheet view ---
$ws1 = $wb.Worksheets.Item("BP4D")

$ws1.Range("A20").Value = "BP4D trained (static) logistic regression"
$row20 = @(0.3674,0.4439,0.3716,0.2929,0.3904,0.3341,0.3796,0.5258,0.4402,0.7204,0.7918,0.7529,0.6669,0.8536,0.7385,0.7156,0.9489,0.8158,0.8674,0.8187,0.8414,0.5464,0.6949,0.6036,0.3149,0.6281,0.4184,0.5148,0.7424,0.5931,0.4381,0.4296,0.4218)
$ws1.Range("B20:AH20").Value = $row20

$ws1.Range("A21").Value = "BP4D trained (static) logistic regression - tanh"
$row21 = @(0.3196,0.5007,0.3891,0.3054,0.4216,0.3542,0.3878,0.5534,0.4559,0.6592,0.8257,0.7331,0.6816,0.8784,0.7676,0.7739,0.8873,0.8266,0.8401,0.882,0.8605,0.5218,0.804,0.6325,0.314,0.6841,0.4292,0.4797,0.8104,0.6026,0.3567,0.5688,0.4359)
$ws1.Range("B21:AH21").Value = $row21

$ws1.Range("AI20").Formula = "=AVERAGE(B20,E20,H20,K20,N20,Q20,T20,W20,Z20,AC20,AF20)"
$ws1.Range("AJ20").Formula = "=AVERAGE(C20,F20,I20,L20,O20,R20,U20,X20,AA20,AD20,AG20)"
$ws1.Range("AK20").Formula = "=AVERAGE(D20,G20,J20,M20,P20,S20,V20,Y20,AB20,AE20,AH20)"

$ws1.Range("AI21").Formula = "=AVERAGE(B21,E21,H21,K21,N21,Q21,T21,W21,Z21,AC21,AF21)"
$ws1.Range("AJ21").Formula = "=AVERAGE(C21,F21,I21,L21,O21,R21,U21,X21,AA21,AD21,AG21)"
$ws1.Range("AK21").Formula = "=AVERAGE(D21,G21,J21,M21,P21,S21,V21,Y21,AB21,AE21,AH21)"

# Rows 22-27: full B:AH number format
$ws1.Range("B22:AH27").NumberFormat = "0.00"
# Row 28: B:Q and V:AH (R:U skipped)
$ws1.Range("B28:Q28").NumberFormat = "0.00"
$ws1.Range("V28:AH28").NumberFormat = "0.00"

$ws1.Range("A23").Select()

$wb.Worksheets.Item("BP4D_intensity").Range("E19").Select()
